$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new values for columns D, M, N, O, P, S
# (D = Fecha, M = Volumen, N = Precio minimo, O = Precio maximo, P = Precio promedio ponderado, S = Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44400; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 3;  D = 44382; M = 24; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 4;  D = 44385; M = 36; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 5;  D = 44291; M = 70; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ Row = 6;  D = 44305; M = 20; N = 22000; O = 22000; P = 22000; S = 1100 },
    @{ Row = 7;  D = 44294; M = 25; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ Row = 8;  D = 44377; M = 25; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 9;  D = 44403; M = 50; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 10; D = 44445; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 11; D = 44448; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 },
    @{ Row = 12; D = 44389; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 13; D = 44300; M = 45; N = 22000; O = 22000; P = 22000; S = 1100 },
    @{ Row = 14; D = 44406; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 15; D = 44301; M = 38; N = 22000; O = 22000; P = 22000; S = 1100 },
    @{ Row = 16; D = 44298; M = 65; N = 22000; O = 22000; P = 22000; S = 1100 },
    @{ Row = 17; D = 44292; M = 30; N = 25000; O = 25000; P = 25000; S = 1250 },
    @{ Row = 18; D = 44307; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 },
    @{ Row = 19; D = 44413; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 },
    @{ Row = 20; D = 44376; M = 38; N = 20000; O = 20000; P = 20000; S = 1000 }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value  = $entry.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $entry.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $entry.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $entry.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $entry.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $entry.S   # S: Precio $/Kg
}
